# HTH: Chapter 15, fix
# "Gayoon sat silent of the chair" -> "Gayoon sat silent on the chair"

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$result = $find.Execute(
    "Gayoon sat silent of the chair",  # FindText
    $true,                              # MatchCase
    $false,                             # MatchWholeWord
    $false,                             # MatchWildcards
    $false,                             # MatchSoundsLike
    $false,                             # MatchAllWordForms
    $true,                              # Forward
    1,                                  # Wrap (wdFindContinue)
    $false,                             # Format
    "Gayoon sat silent on the chair",  # ReplaceWith
    2                                   # Replace (wdReplaceAll)
)

if (-not $result) {
    throw "Could not find target text 'Gayoon sat silent of the chair' to replace."
}

Write-Host "Replaced 'sat silent of the chair' -> 'sat silent on the chair'"
